$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1697.8
$ws.Range("I33").Value = 1544.4419
$ws.Range("K33").Value = 1544.4419
$ws.Range("M33").Value = -1315.4419
$ws.Range("H43").Value = 17665.166
$ws.Range("J43").Value = 10549
$ws.Range("L43").Value = 10549
$ws.Range("N43").Value = -10687
$ws.Range("H70").Value = 8217.799999999999
$ws.Range("I70").Value = 5832.3335
$ws.Range("J70").Value = 9240.143
$ws.Range("K70").Value = 17497.0005
$ws.Range("L70").Value = 27720.429
$ws.Range("M70").Value = -17227.0005
$ws.Range("N70").Value = -28260.429
$ws.Range("H73").Value = 8217.799999999999
$ws.Range("I73").Value = 5832.3335
$ws.Range("J73").Value = 9240.143
$ws.Range("K73").Value = 17497.0005
$ws.Range("L73").Value = 27720.429
$ws.Range("M73").Value = -16561.0005
$ws.Range("N73").Value = -29592.429
$ws.Range("H138").Value = 3366
$ws.Range("I138").Value = 1512.2307
$ws.Range("K138").Value = 4536.6921
$ws.Range("M138").Value = 603.3078999999998
$ws.Range("H141").Value = 3104.3333
$ws.Range("J141").Value = 5329.6665
$ws.Range("L141").Value = 15988.9995
$ws.Range("N141").Value = -26348.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1282.1111
$ws.Range("I32").Value = 865.1475
$ws.Range("K32").Value = 865.1475
$ws.Range("M32").Value = -578.1475
$ws.Range("H45").Value = 7856.303
$ws.Range("J45").Value = 1120.0769
$ws.Range("L45").Value = 1120.0769
$ws.Range("N45").Value = -1874.0769
$ws.Range("H74").Value = 1990
$ws.Range("I74").Value = 1808.5
$ws.Range("K74").Value = 1808.5
$ws.Range("M74").Value = -934.5
$ws.Range("H77").Value = 1990
$ws.Range("I77").Value = 1808.5
$ws.Range("K77").Value = 9042.5
$ws.Range("M77").Value = -4674.5
$ws.Range("H102").Value = 2496.9167
$ws.Range("I102").Value = 2587.5454
$ws.Range("K102").Value = 2587.5454
$ws.Range("M102").Value = -965.5454
$ws.Range("H132").Value = 2926.8928
$ws.Range("I132").Value = 2748.875
$ws.Range("K132").Value = 8246.625
$ws.Range("M132").Value = -5716.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3892.75
$ws.Range("I99").Value = 3671.3
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 3671.3
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -2173.3
$ws.Range("N99").Value = -7996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 90000
$ws.Range("J20").Value = 90000
$ws.Range("L20").Value = 90000
$ws.Range("N20").Value = -90472
$ws.Range("H30").Value = 90000
$ws.Range("J30").Value = 90000
$ws.Range("L30").Value = 90000
$ws.Range("N30").Value = -90182
$ws.Range("H31").Value = 4724.6787
$ws.Range("J31").Value = 3640.889
$ws.Range("L31").Value = 3640.889
$ws.Range("N31").Value = -4230.889
$ws.Range("H34").Value = 4724.6787
$ws.Range("J34").Value = 3640.889
$ws.Range("L34").Value = 3640.889
$ws.Range("N34").Value = -4044.889
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H99").Value = 5274
$ws.Range("I99").Value = 3699.3333
$ws.Range("K99").Value = 3699.3333
$ws.Range("M99").Value = -2201.3333
$ws.Range("H126").Value = 5274
$ws.Range("I126").Value = 3699.3333
$ws.Range("K126").Value = 11097.9999
$ws.Range("M126").Value = -8627.999899999999
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960
$ws.Range("H132").Value = 4210.6562
$ws.Range("I132").Value = 4023.9355
$ws.Range("K132").Value = 12071.8065
$ws.Range("M132").Value = -9541.806500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69997
$ws.Range("J37").Value = 69997
$ws.Range("L37").Value = 209991
$ws.Range("N37").Value = -210215
$ws.Range("H75").Value = 2717.7856
$ws.Range("J75").Value = 2727.2222
$ws.Range("L75").Value = 8181.6666
$ws.Range("N75").Value = -10177.6666
$ws.Range("H78").Value = 2717.7856
$ws.Range("J78").Value = 2727.2222
$ws.Range("L78").Value = 24544.9998
$ws.Range("N78").Value = -34528.99980000001
$ws.Range("H128").Value = 130000
$ws.Range("I128").Value = 130000
$ws.Range("K128").Value = 390000
$ws.Range("M128").Value = -385020
$ws.Range("H136").Value = 14009.429
$ws.Range("I136").Value = 14453.2
$ws.Range("J136").Value = 12900
$ws.Range("K136").Value = 43359.60000000001
$ws.Range("L136").Value = 38700
$ws.Range("M136").Value = -38259.60000000001
$ws.Range("N136").Value = -48900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1927.25
$ws.Range("J80").Value = 1952.5
$ws.Range("L80").Value = 1952.5
$ws.Range("N80").Value = -3948.5
$ws.Range("H83").Value = 1927.25
$ws.Range("J83").Value = 1952.5
$ws.Range("L83").Value = 9762.5
$ws.Range("N83").Value = -19746.5
$ws.Range("H132").Value = 5177.9585
$ws.Range("I132").Value = 4739.636
$ws.Range("K132").Value = 14218.908
$ws.Range("M132").Value = -11688.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3408.4783
$ws.Range("I100").Value = 2103.9167
$ws.Range("J100").Value = 4831.636
$ws.Range("K100").Value = 2103.9167
$ws.Range("L100").Value = 4831.636
$ws.Range("M100").Value = -1562.9167
$ws.Range("N100").Value = -5913.636
$ws.Range("H132").Value = 5405.3125
$ws.Range("I132").Value = 3868.7144
$ws.Range("J132").Value = 6600.4443
$ws.Range("K132").Value = 11606.1432
$ws.Range("L132").Value = 19801.3329
$ws.Range("M132").Value = -9076.143199999999
$ws.Range("N132").Value = -24861.3329
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450
$ws.Range("H138").Value = 81628.57000000001
$ws.Range("J138").Value = 81628.57000000001
$ws.Range("L138").Value = 81628.57000000001
$ws.Range("N138").Value = -91908.57000000001
$ws.Range("H139").Value = 82677.75
$ws.Range("J139").Value = 82677.75
$ws.Range("L139").Value = 82677.75
$ws.Range("N139").Value = -92957.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2995.6765
$ws.Range("I132").Value = 3328.4783
$ws.Range("J132").Value = 2299.818
$ws.Range("K132").Value = 9985.4349
$ws.Range("L132").Value = 6899.454000000001
$ws.Range("M132").Value = -7455.4349
$ws.Range("N132").Value = -11959.454
$ws.Range("H136").Value = 2615.1667
$ws.Range("I136").Value = 2140.52
$ws.Range("J136").Value = 3693.9092
$ws.Range("K136").Value = 6421.559999999999
$ws.Range("L136").Value = 11081.7276
$ws.Range("M136").Value = -3871.559999999999
$ws.Range("N136").Value = -16181.7276
